$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.800999999999999
$ws.Range("B8").Value = 6.185
$ws.Range("B10").Value = 5.777
$ws.Range("B12").Value = 5.090000000000001
$ws.Range("D15").Value = -8.147000000000002
$ws.Range("B18").Value = 5.040999999999999
$ws.Range("D18").Value = -8.320000000000002
$ws.Range("D20").Value = -7.702000000000001
$ws.Range("D29").Value = -7.476999999999999
$ws.Range("D30").Value = -7.179
$ws.Range("D31").Value = -8.086
$ws.Range("B37").Value = 8.73
$ws.Range("D40").Value = -7.663000000000001
$ws.Range("D50").Value = -8.125999999999999
$ws.Range("B55").Value = 4.621
$ws.Range("B68").Value = 5.100999999999999
$ws.Range("D68").Value = -6.778
$ws.Range("D76").Value = -7.444
$ws.Range("B77").Value = 5.880000000000001
$ws.Range("B78").Value = 7.811
$ws.Range("B81").Value = 6.058000000000001
$ws.Range("B82").Value = 5.659000000000001
$ws.Range("D87").Value = -8.261999999999999
$ws.Range("D88").Value = -8.09
$ws.Range("D96").Value = -7.267
$ws.Range("D98").Value = -8.228000000000002
$ws.Range("D101").Value = -7.886000000000001
$ws.Range("D102").Value = -8.036
